# Weekly Fruit/Vegetable price update — "Perejil" (Parsley) sheet.
# A new weekly observation is inserted as row 201 (pushing the existing
# rows 201-228 down to 202-229), growing the used range from A1:R228 to
# A1:R229.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 201; everything currently at/after
# row 201 shifts down by one (old row 201 -> 202, ..., old row 228 -> 229).
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row 201 with the new weekly record.
$ws.Cells.Item(201, 1).Value = 9
$ws.Cells.Item(201, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(201, 3).Value = "Metropolitana"
$ws.Cells.Item(201, 4).Value = 44491
$ws.Cells.Item(201, 5).Value = 13
$ws.Cells.Item(201, 6).Value = 100112044
$ws.Cells.Item(201, 7).Value = "Perejil"
$ws.Cells.Item(201, 8).Value = "Sin especificar"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 106
$ws.Cells.Item(201, 11).Value = 9000
$ws.Cells.Item(201, 12).Value = 10000
$ws.Cells.Item(201, 13).Value = 9500
$ws.Cells.Item(201, 14).Value = "$/docena de atados"
$ws.Cells.Item(201, 15).Value = "Región Metropolitana"
$ws.Cells.Item(201, 16).Value = 3167
$ws.Cells.Item(201, 17).Value = 3
$ws.Cells.Item(201, 18).Value = "Hortaliza"
